# Apply updated cryptocurrency market data (price + 1h volume change)
# Rows 2-51 of Sheet1; some rows also swap Coin/Link (re-ranked entries).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value.
# Column D (Price) values that Excel would otherwise auto-convert to a
# number (dropping formatting like trailing zeros, e.g. "1.00" -> 1) are
# prefixed with a leading apostrophe to force them to stay text, exactly
# as they were already stored (t="inlineStr") in the workbook.
$updates = [ordered]@{
    'D2' = '68.192.39'
    'E2' = '  +1.88%  '
    'D3' = '2.628.37'
    'E3' = '  +1.20%  '
    'D4' = '''0.999'
    'E4' = '  -0.28%  '
    'D5' = '''597.47'
    'E5' = '  +0.98%  '
    'D6' = '''153.95'
    'E6' = '  +1.42%  '
    'E7' = '  -0.10%  '
    'D8' = '''0.545'
    'E8' = '  -0.97%  '
    'D9' = '2.626.00'
    'E9' = '  +1.17%  '
    'E10' = '  +10.93%  '
    'E11' = '  -0.51%  '
    'D12' = '''5.23'
    'E12' = '  +0.77%  '
    'E13' = '  +0.50%  '
    'D14' = '''27.61'
    'E14' = '  +0.32%  '
    'D15' = '''0.0000188'
    'E15' = '  +5.08%  '
    'D16' = '3.094.81'
    'E16' = '  +0.58%  '
    'D17' = '68.012.34'
    'E17' = '  +1.78%  '
    'D18' = '2.621.18'
    'E18' = '  +1.29%  '
    'D19' = '''11.41'
    'E19' = '  +3.96%  '
    'D20' = '''368.21'
    'E20' = '  +0.61%  '
    'D21' = '''7.43'
    'E21' = '  +1.07%  '
    'D22' = '''4.23'
    'E22' = '  -1.27%  '
    'D23' = '''4.82'
    'E23' = '  -0.27%  '
    'D24' = '''2.08'
    'E24' = '  +1.50%  '
    'D25' = '''72.71'
    'E25' = '  +7.79%  '
    'E26' = '  +0.05%  '
    'D27' = '''9.96'
    'E27' = '  -1.23%  '
    'B28' = 'PEPE'
    'C28' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D28' = '''0.0000104'
    'E28' = '  +3.93%  '
    'B29' = 'WrappedeETH'
    'C29' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D29' = '2.750.65'
    'E29' = '  +0.41%  '
    'D30' = '''1.00'
    'E30' = '  -0.09%  '
    'D31' = '''573.68'
    'E31' = '  -1.50%  '
    'B32' = 'Fetch.AI'
    'C32' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D32' = '''1.40'
    'E32' = '  +1.10%  '
    'B33' = 'InternetComputer(DFINITY)'
    'C33' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D33' = '''7.89'
    'E33' = '  +2.49%  '
    'D34' = '''1.84'
    'E34' = '  +1.98%  '
    'D35' = '''0.999'
    'E35' = '  -0.09%  '
    'E36' = '  +3.63%  '
    'E37' = '  +1.85%  '
    'D38' = '''162.28'
    'E38' = '  +5.05%  '
    'D39' = '''19.13'
    'E39' = '  +1.58%  '
    'D40' = '''1.89'
    'E40' = '  +5.27%  '
    'D41' = '''0.367'
    'E41' = '  +0.71%  '
    'D42' = '''5.32'
    'E42' = '  +2.20%  '
    'D43' = '''2.66'
    'E43' = '  +3.92%  '
    'D44' = '''17.55'
    'E44' = '  +4.54%  '
    'D45' = '0.0₆0323'
    'E45' = '  +10.46%  '
    'E46' = '  +0.08%  '
    'D47' = '''40.20'
    'E47' = '  -1.45%  '
    'D48' = '''155.01'
    'E48' = '  +0.83%  '
    'D49' = '''3.69'
    'E49' = '  -0.53%  '
    'D50' = '''21.89'
    'E50' = '  +1.84%  '
    'E51' = '  +0.35%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
